$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.950.30'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '3.725.97'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.54'
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.26'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("D7").Value = '3.721.76'
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  +4.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.23'
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.26'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '4.341.51'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '3.719.01'
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").Value = '67.892.07'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.33'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.43'
$ws.Range("E20").Value = '  +8.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.97'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.03'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000143'
$ws.Range("E25").Value = '  +4.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.36'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.38'
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.79'
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.58'
$ws.Range("E33").Value = '  -2.91%  '
$ws.Range("D34").Value = '3.857.61'
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").Value = '3.662.79'
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.87'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '48.98'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '432.72'
$ws.Range("E43").Value = '  -3.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.97'
$ws.Range("E44").Value = '  -2.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.47'
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.94'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.26'
$ws.Range("E49").Value = '  +2.15%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = '2.764.11'
$ws.Range("E51").Value = '  -3.32%  '
